$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Migrate the "level" column (D) to sit right after "id" (A), i.e. into
# column B, shifting the former storeLimit/priceMultiplier columns (B, C)
# one slot to the right (into C, D). This reflects the
# FarmerStatTable -> FarmerLevelTable column cleanup: level now leads the
# stat columns. Swap the cell values row-by-row (header, type row and all
# data rows) rather than doing a structural column insert, so existing
# per-column formatting (styles) stays attached to the same column letters
# as before.
for ($r = 2; $r -le 17; $r++) {
    $levelVal = $ws.Cells.Item($r, 4).Value2
    $storeLimitVal = $ws.Cells.Item($r, 2).Value2
    $priceMultVal = $ws.Cells.Item($r, 3).Value2

    $ws.Cells.Item($r, 2).Value2 = $levelVal
    $ws.Cells.Item($r, 3).Value2 = $storeLimitVal
    $ws.Cells.Item($r, 4).Value2 = $priceMultVal
}

# --- Freeze panes so the id/level columns and header rows stay pinned while
# scrolling through the table, then leave the viewport scrolled to around
# column P, row 5 (matches the saved view state).
$ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("P5").Select()
